$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-19 19:48:27"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "69%"
$ws.Range("I2").Value = "2.9 mm"
$ws.Range("E3").Value = "2026-02-19 19:48:29"
$ws.Range("I3").Value = "4.4 mm"
$ws.Range("E4").Value = "2026-02-19 19:48:32"
$ws.Range("J4").Value = "1009.7 hPa"
$ws.Range("E5").Value = "2026-02-19 19:48:34"
$ws.Range("I5").Value = "7.4 mm"
$ws.Range("E6").Value = "2026-02-19 19:48:37"
$ws.Range("J6").Value = "1009.8 hPa"
$ws.Range("E7").Value = "2026-02-19 19:48:39"
$ws.Range("J7").Value = "1010.8 hPa"
$ws.Range("O7").Value = "13.9 °C"
$ws.Range("E8").Value = "2026-02-19 19:48:42"
$ws.Range("J8").Value = "1010.4 hPa"
$ws.Range("E9").Value = "2026-02-19 19:48:44"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "78%"
$ws.Range("E10").Value = "2026-02-19 19:48:47"
$ws.Range("N10").Value = "5.0 °C 19:29 TU"
$ws.Range("O10").Value = "10.7 °C"
$ws.Range("E11").Value = "2026-02-19 19:48:49"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "65%"
$ws.Range("O11").Value = "5.6 °C"
$ws.Range("E12").Value = "2026-02-19 19:48:51"
$ws.Range("E13").Value = "2026-02-19 19:48:54"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "63%"
$ws.Range("J13").Value = "1011.0 hPa"
$ws.Range("E14").Value = "2026-02-19 19:48:56"
$ws.Range("O14").Value = "13.3 °C"
$ws.Range("E15").Value = "2026-02-19 19:48:59"
$ws.Range("O15").Value = "10.0 °C"
$ws.Range("E16").Value = "2026-02-19 19:49:01"
$ws.Range("I16").Value = "8.9 mm"
$ws.Range("E17").Value = "2026-02-19 19:49:03"
$ws.Range("E18").Value = "2026-02-19 19:49:06"
$ws.Range("J18").Value = "1010.0 hPa"
$ws.Range("E19").Value = "2026-02-19 19:49:08"
$ws.Range("E20").Value = "2026-02-19 19:49:11"
$ws.Range("L20").Value = "88.2 km/h - 341º 19:22 TU"
$ws.Range("E21").Value = "2026-02-19 19:49:13"
$ws.Range("J21").Value = "1011.0 hPa"
$ws.Range("E22").Value = "2026-02-19 19:49:16"
$ws.Range("E23").Value = "2026-02-19 19:49:18"
$ws.Range("G23").Value = "215 cm"
$ws.Range("I23").Value = "9.3 mm"
$ws.Range("E24").Value = "2026-02-19 19:49:20"
$ws.Range("J24").Value = "1014.6 hPa"
$ws.Range("O24").Value = "9.0 °C"
$ws.Range("E25").Value = "2026-02-19 19:49:23"
$ws.Range("I25").Value = "5.3 mm"
$ws.Range("O25").Value = "-4.4 °C"
$ws.Range("E26").Value = "2026-02-19 19:49:25"
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = "57%"
$ws.Range("J26").Value = "1009.7 hPa"
$ws.Range("E27").Value = "2026-02-19 19:49:28"
$ws.Range("O27").Value = "-3.8 °C"
$ws.Range("E28").Value = "2026-02-19 19:49:30"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "66%"
$ws.Range("J28").Value = "1009.6 hPa"
$ws.Range("L28").Value = "33.1 km/h - 283º 19:17 TU"
$ws.Range("O28").Value = "9.3 °C"
$ws.Range("E29").Value = "2026-02-19 19:49:33"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "75%"
$ws.Range("O29").Value = "10.8 °C"
$ws.Range("E30").Value = "2026-02-19 19:49:35"
$ws.Range("J30").Value = "1009.8 hPa"
$ws.Range("E31").Value = "2026-02-19 19:49:38"
$ws.Range("J31").Value = "1009.3 hPa"
$ws.Range("E32").Value = "2026-02-19 19:49:40"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "68%"
$ws.Range("O32").Value = "4.8 °C"
$ws.Range("E33").Value = "2026-02-19 19:49:43"
$ws.Range("J33").Value = "1010.6 hPa"
$ws.Range("E34").Value = "2026-02-19 19:49:45"
$ws.Range("E35").Value = "2026-02-19 19:49:48"
$ws.Range("J35").Value = "1016.0 hPa"
$ws.Range("E36").Value = "2026-02-19 19:49:50"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "65%"
$ws.Range("J36").Value = "1010.1 hPa"
$ws.Range("E37").Value = "2026-02-19 19:49:53"
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = "69%"
$ws.Range("J37").Value = "1011.1 hPa"
$ws.Range("E38").Value = "2026-02-19 19:49:55"
$ws.Range("E39").Value = "2026-02-19 19:49:58"
$ws.Range("I39").Value = "4.6 mm"
$ws.Range("E40").Value = "2026-02-19 19:50:00"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "72%"
$ws.Range("J40").Value = "1012.2 hPa"
$ws.Range("O40").Value = "6.4 °C"
$ws.Range("E41").Value = "2026-02-19 19:50:03"
$ws.Range("J41").Value = "1012.7 hPa"
$ws.Range("E42").Value = "2026-02-19 19:50:05"
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = "75%"
$ws.Range("E43").Value = "2026-02-19 19:50:07"
$ws.Range("E44").Value = "2026-02-19 19:50:10"
$ws.Range("I44").Value = "8.8 mm"
$ws.Range("E45").Value = "2026-02-19 19:50:12"
$ws.Range("J45").Value = "1015.2 hPa"
$ws.Range("E46").Value = "2026-02-19 19:50:15"
$ws.Range("J46").Value = "1015.4 hPa"
